$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.210.43"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "1.917.03"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  -0.99%  "
$ws.Range("D5").Value = "317.36"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").Value = "0.3832"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "0.07381"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "0.9420"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "20.96"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "0.07816"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "1.904.62"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "5.506"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "6.650"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "91.39"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D18").Value = "0.000008847"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "28.222.99"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("D21").Value = "14.88"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "5.163"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "2.164.06"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").Value = "156.35"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "1.925"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "2.104"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("D29").Value = "116.41"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "4.979"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "0.08913"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "3.357"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "1.252"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "0.7753"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("D35").Value = "4.707"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("D36").Value = "2.687"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "0.02054"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "1.105"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "0.5545"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("D40").Value = "0.05334"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "3.005"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "7.055"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "0.1533"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "8.484"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "10.78"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("D46").Value = "0.4869"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "107.01"
$ws.Range("E47").Value = "  +3.76%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "1.661"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "68.66"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").Value = "0.06123"
$ws.Range("E51").Value = "  +0.14%  "
